$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 423, shifting existing rows 423-444 down to 424-445
$ws.Rows.Item(423).Insert()

# Fill in the new row 423 with the data from the diff
$ws.Range("A423").Value = 5
$ws.Range("B423").Value = "Macroferia Regional de Talca"
$ws.Range("C423").Value = "Maule"
$ws.Range("D423").Value = 44753
$ws.Range("E423").Value = 7
$ws.Range("F423").Value = 100112043
$ws.Range("G423").Value = "Pepino ensalada"
$ws.Range("H423").Value = "Sin especificar"
$ws.Range("I423").Value = "Primera"
$ws.Range("J423").Value = 300
$ws.Range("K423").Value = 18000
$ws.Range("L423").Value = 18000
$ws.Range("M423").Value = 18000
$ws.Range("N423").Value = "$/caja 60 unidades"
$ws.Range("O423").Value = "Región de Arica y Parinacota"
$ws.Range("P423").Value = 300
$ws.Range("Q423").Value = 60
$ws.Range("R423").Value = "Hortaliza"
